# Weekly update: insert a new price record as row 21, pushing the
# previously existing rows 21-67 down to 22-68 (dimension grows to R68).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 21 (shifts rows 21:67 -> 22:68,
# carrying formatting/styles along, same as a native Excel row insert).
$ws.Rows.Item(21).Insert()

# Populate the newly inserted row 21 with this week's record.
$ws.Range("A21").Value = 3
$ws.Range("B21").Value = "Femacal de La Calera"
$ws.Range("C21").Value = "Coquimbo"
$ws.Range("D21").Value = 44868
$ws.Range("E21").Value = 5
$ws.Range("F21").Value = 100112022
$ws.Range("G21").Value = "Arveja Verde"
$ws.Range("H21").Value = "Perfection"
$ws.Range("I21").Value = "Primera"
$ws.Range("J21").Value = 76
$ws.Range("K21").Value = 22000
$ws.Range("L21").Value = 23000
$ws.Range("M21").Value = 22500
$ws.Range("N21").Value = "$/malla 25 kilos"
$ws.Range("O21").Value = "Provincia de Quillota"
$ws.Range("P21").Value = 900
$ws.Range("Q21").Value = 25
$ws.Range("R21").Value = "Hortaliza"
